$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.577.78"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.955.86"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.29"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.31"
$ws.Range("E7").Value = "  +1.81%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +3.53%  "

$ws.Range("E10").Value = "  -6.21%  "

$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.04"
$ws.Range("E12").Value = "  +4.11%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.838"
$ws.Range("E13").Value = "  +2.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.241.68"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.25"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.29"
$ws.Range("E16").Value = "  +1.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.950.95"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.455.02"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.74"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.67"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.05"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +2.85%  "

$ws.Range("E25").Value = "  +2.84%  "

$ws.Range("E26").Value = "  +6.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.23"
$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.31"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +5.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.73"
$ws.Range("E32").Value = "  +2.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0611"
$ws.Range("E33").Value = "  -3.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.41"
$ws.Range("E34").Value = "  +4.75%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.47"
$ws.Range("E35").Value = "  +14.16%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.28"
$ws.Range("E36").Value = "  +6.42%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").Value = "  -12.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0980"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("E41").Value = "  +1.47%  "

$ws.Range("E42").Value = "  +0.19%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.370.45"
$ws.Range("E44").Value = "  +2.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.70"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.92"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.12"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.132.00"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.05"
$ws.Range("E51").Value = "  -1.21%  "
